# Applies the CodeSystem-attribution-product.xlsx update:
#  - Metadata sheet: remove the duplicate "Contact" row, replace the
#    "Contact" row with "Jurisdiction" / "United States of America",
#    populate "Publisher" value with "Alvearie Team", populate
#    "Case Sensitive" value with "true", bump Version to 6.0.0 and
#    Date to 2022-01-21T20:46:54+00:00.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "6.0.0"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Fill in the Publisher value (row 9, column B), which was previously blank
$ws.Range("B9").Value = "Alvearie Team"

# Delete the second, duplicate "Contact" row (row 11) entirely, which
# shifts all subsequent rows up by one.
$ws.Rows.Item(11).Delete()

# The remaining "Contact" row (now row 10) becomes "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# "Case Sensitive" row (now row 14 after the shift) gets a text value of
# "true" (not the Boolean TRUE). A leading apostrophe forces Excel to
# store it as literal text, matching the shared-string based target.
$ws.Range("B14").Value = "'true"
